$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.897.70'
$ws.Range('E2').Value = '  +0.63%  '
$ws.Range('D3').Value = '1.840.15'
$ws.Range('E3').Value = '  +1.95%  '
$ws.Range('E4').Value = '  +0.24%  '
$ws.Range('D5').Value = "'232.24"
$ws.Range('E5').Value = '  +0.64%  '
$ws.Range('E6').Value = '  +2.72%  '
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('D8').Value = "'41.09"
$ws.Range('E8').Value = '  +5.83%  '
$ws.Range('E9').Value = '  +3.69%  '
$ws.Range('D10').Value = "'0.0690"
$ws.Range('E10').Value = '  +2.02%  '
$ws.Range('E11').Value = '  -0.87%  '
$ws.Range('D12').Value = '2.108.61'
$ws.Range('E12').Value = '  +2.08%  '
$ws.Range('D13').Value = "'11.35"
$ws.Range('E13').Value = '  +4.72%  '
$ws.Range('D14').Value = '1.840.26'
$ws.Range('E14').Value = '  +3.14%  '
$ws.Range('D15').Value = "'0.670"
$ws.Range('E15').Value = '  +1.87%  '
$ws.Range('D16').Value = "'4.65"
$ws.Range('E16').Value = '  +2.70%  '
$ws.Range('D17').Value = '34.916.18'
$ws.Range('E17').Value = '  +0.74%  '
$ws.Range('D18').Value = "'69.82"
$ws.Range('E18').Value = '  +0.68%  '
$ws.Range('D19').Value = '0.0₃0788'
$ws.Range('E19').Value = '  +1.19%  '
$ws.Range('D20').Value = "'239.96"
$ws.Range('E20').Value = '  +0.47%  '
$ws.Range('E21').Value = '  +3.79%  '
$ws.Range('D22').Value = "'4.75"
$ws.Range('E22').Value = '  +2.64%  '
$ws.Range('E23').Value = '  +0.13%  '
$ws.Range('D24').Value = "'2.25"
$ws.Range('E24').Value = '  +1.30%  '
$ws.Range('D25').Value = "'171.95"
$ws.Range('E25').Value = '  -0.28%  '
$ws.Range('D26').Value = "'7.85"
$ws.Range('E26').Value = '  +2.14%  '
$ws.Range('D27').Value = "'17.39"
$ws.Range('E27').Value = '  +1.92%  '
$ws.Range('E28').Value = '  +3.94%  '
$ws.Range('D29').Value = "'1.64"
$ws.Range('E29').Value = '  +8.73%  '
$ws.Range('E30').Value = '  +0.29%  '
$ws.Range('E31').Value = '  +1.73%  '
$ws.Range('D32').Value = "'3.95"
$ws.Range('E32').Value = '  -0.75%  '
$ws.Range('E33').Value = '  -0.80%  '
$ws.Range('D34').Value = "'1.62"
$ws.Range('E34').Value = '  +22.31%  '
$ws.Range('D35').Value = "'1.94"
$ws.Range('E35').Value = '  +10.57%  '
$ws.Range('D36').Value = "'1.25"
$ws.Range('E36').Value = '  +0.26%  '
$ws.Range('D37').Value = "'0.744"
$ws.Range('E37').Value = '  +8.41%  '
$ws.Range('E38').Value = '  +11.06%  '
$ws.Range('D39').Value = "'89.70"
$ws.Range('E39').Value = '  -0.89%  '
$ws.Range('E40').Value = '  +3.72%  '
$ws.Range('D41').Value = '1.337.35'
$ws.Range('E41').Value = '  +2.44%  '
$ws.Range('D42').Value = "'14.50"
$ws.Range('E42').Value = '  +2.46%  '
$ws.Range('E43').Value = '  -1.91%  '
$ws.Range('E44').Value = '  +1.98%  '
$ws.Range('E45').Value = '  +3.42%  '
$ws.Range('E46').Value = '  +3.68%  '
$ws.Range('E47').Value = '  +2.99%  '
$ws.Range('B48').Value = 'RocketPoolETH'
$ws.Range('C48').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D48').Value = '2.025.50'
$ws.Range('E48').Value = '  +1.62%  '
$ws.Range('B49').Value = 'Gas'
$ws.Range('C49').Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range('D49').Value = "'10.70"
$ws.Range('E49').Value = '  +66.12%  '
$ws.Range('E50').Value = '  +0.16%  '
$ws.Range('D51').Value = "'3.38"
$ws.Range('E51').Value = '  +15.91%  '
